# Update timestamps on the "Generate Report for Handback" run.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-04 17:10:36"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-04 17:10:32"
$wsZhCn.Range("K2").Value = "2016-09-04 17:10:57"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-04 17:10:36"
$wsDeDe.Range("K2").Value = "2016-09-04 17:11:10"
